$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44522
$ws.Range("M2").Value = 25
$ws.Range("N2").Value = 30000
$ws.Range("O2").Value = 30000
$ws.Range("P2").Value = 30000
$ws.Range("S2").Value = 1500

# Row 3
$ws.Range("D3").Value = 44417
$ws.Range("N3").Value = 24000
$ws.Range("O3").Value = 24000
$ws.Range("P3").Value = 24000
$ws.Range("S3").Value = 1200

# Row 4
$ws.Range("D4").Value = 44305
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = 24000
$ws.Range("O4").Value = 24000
$ws.Range("P4").Value = 24000
$ws.Range("S4").Value = 1200

# Row 5
$ws.Range("D5").Value = 44372
$ws.Range("M5").Value = 60
$ws.Range("P5").Value = 20667
$ws.Range("S5").Value = 1033

# Row 6
$ws.Range("D6").Value = 44435
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 25000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 25000
$ws.Range("S6").Value = 1250

# Row 7
$ws.Range("D7").Value = 44445
$ws.Range("M7").Value = 35
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 20000
$ws.Range("S7").Value = 1000

# Row 8
$ws.Range("D8").Value = 44356
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 21000
$ws.Range("P8").Value = 20500
$ws.Range("S8").Value = 1025

# Row 9
$ws.Range("D9").Value = 44438
$ws.Range("M9").Value = 25
$ws.Range("N9").Value = 21000
$ws.Range("P9").Value = 21000
$ws.Range("S9").Value = 1050

# Row 10
$ws.Range("D10").Value = 44333
$ws.Range("M10").Value = 30

# Row 11
$ws.Range("D11").Value = 44410
$ws.Range("N11").Value = 25000
$ws.Range("O11").Value = 25000
$ws.Range("P11").Value = 25000
$ws.Range("S11").Value = 1250

# Row 12
$ws.Range("D12").Value = 44355
$ws.Range("M12").Value = 200
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 21000
$ws.Range("P12").Value = 20500
$ws.Range("R12").Value = "Ecuador"
$ws.Range("S12").Value = 1025

# Row 13
$ws.Range("D13").Value = 44442
$ws.Range("M13").Value = 30
$ws.Range("N13").Value = 22000
$ws.Range("O13").Value = 22000
$ws.Range("P13").Value = 22000
$ws.Range("S13").Value = 1100

# Row 14
$ws.Range("D14").Value = 44363
$ws.Range("M14").Value = 150
$ws.Range("N14").Value = 21000
$ws.Range("P14").Value = 21500
$ws.Range("S14").Value = 1075

# Row 15
$ws.Range("D15").Value = 44284
$ws.Range("M15").Value = 40
$ws.Range("N15").Value = 23000
$ws.Range("O15").Value = 23000
$ws.Range("P15").Value = 23000
$ws.Range("S15").Value = 1150

# Row 16
$ws.Range("D16").Value = 44300
$ws.Range("N16").Value = 19000
$ws.Range("O16").Value = 20000
$ws.Range("P16").Value = 19500
$ws.Range("S16").Value = 975

# Row 17
$ws.Range("D17").Value = 44382
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 19000
$ws.Range("O17").Value = 20000
$ws.Range("P17").Value = 19500
$ws.Range("S17").Value = 975

# Row 18
$ws.Range("D18").Value = 44431
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 60
$ws.Range("N18").Value = 25000
$ws.Range("O18").Value = 25000
$ws.Range("P18").Value = 25000
$ws.Range("S18").Value = 1250

# Row 19
$ws.Range("D19").Value = 44424
$ws.Range("M19").Value = 70
$ws.Range("O19").Value = 25000
$ws.Range("P19").Value = 24429
$ws.Range("S19").Value = 1221

# Row 20
$ws.Range("D20").Value = 44326
$ws.Range("M20").Value = 40
$ws.Range("N20").Value = 22000
$ws.Range("O20").Value = 22000
$ws.Range("P20").Value = 22000
$ws.Range("R20").Value = "Perú"
$ws.Range("S20").Value = 1100

# Row 21
$ws.Range("D21").Value = 44354
$ws.Range("M21").Value = 150
$ws.Range("N21").Value = 21000
$ws.Range("O21").Value = 22000
$ws.Range("P21").Value = 21500
$ws.Range("S21").Value = 1075

# Row 22
$ws.Range("D22").Value = 44312
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = 22000
$ws.Range("O22").Value = 22000
$ws.Range("P22").Value = 22000
$ws.Range("S22").Value = 1100

# Row 23
$ws.Range("D23").Value = 44299
$ws.Range("M23").Value = 150
$ws.Range("N23").Value = 19000
$ws.Range("O23").Value = 20000
$ws.Range("P23").Value = 19500
$ws.Range("S23").Value = 975

# Row 24
$ws.Range("D24").Value = 44166
$ws.Range("M24").Value = 120
$ws.Range("N24").Value = 28000
$ws.Range("O24").Value = 28000
$ws.Range("P24").Value = 28000
$ws.Range("S24").Value = 1400

# Row 25
$ws.Range("D25").Value = 44277
$ws.Range("M25").Value = 60
$ws.Range("N25").Value = 24000
$ws.Range("O25").Value = 24000
$ws.Range("P25").Value = 24000
$ws.Range("S25").Value = 1200

# Row 26
$ws.Range("D26").Value = 44165
$ws.Range("M26").Value = 300
$ws.Range("N26").Value = 27000
$ws.Range("O26").Value = 28000
$ws.Range("P26").Value = 27500
$ws.Range("S26").Value = 1375

# Row 27
$ws.Range("D27").Value = 44263
$ws.Range("L27").Value = "Segunda"
$ws.Range("M27").Value = 150
$ws.Range("N27").Value = 15000
$ws.Range("O27").Value = 15000
$ws.Range("P27").Value = 15000
$ws.Range("S27").Value = 750

# Row 28
$ws.Range("D28").Value = 44396
$ws.Range("M28").Value = 45
$ws.Range("N28").Value = 22000
$ws.Range("O28").Value = 22000
$ws.Range("P28").Value = 22000
$ws.Range("S28").Value = 1100

# Row 29
$ws.Range("D29").Value = 44350
$ws.Range("M29").Value = 90
$ws.Range("N29").Value = 21000
$ws.Range("O29").Value = 22000
$ws.Range("P29").Value = 21556
$ws.Range("S29").Value = 1078

# Row 30
$ws.Range("D30").Value = 44529
$ws.Range("M30").Value = 34
$ws.Range("N30").Value = 28000
$ws.Range("O30").Value = 28000
$ws.Range("P30").Value = 28000
$ws.Range("S30").Value = 1400

# Row 31
$ws.Range("D31").Value = 44298
$ws.Range("M31").Value = 240
$ws.Range("N31").Value = 19000
$ws.Range("O31").Value = 20000
$ws.Range("P31").Value = 19500
$ws.Range("S31").Value = 975

# Row 32
$ws.Range("D32").Value = 44452
$ws.Range("M32").Value = 35
$ws.Range("N32").Value = 21000
$ws.Range("O32").Value = 22000
$ws.Range("P32").Value = 21429
$ws.Range("S32").Value = 1071

# Row 34
$ws.Range("D34").Value = 44473
$ws.Range("M34").Value = 40
$ws.Range("N34").Value = 24000
$ws.Range("O34").Value = 24000
$ws.Range("P34").Value = 24000
$ws.Range("S34").Value = 1200

# Row 35
$ws.Range("D35").Value = 44302
$ws.Range("M35").Value = 100
$ws.Range("N35").Value = 19000
$ws.Range("O35").Value = 20000
$ws.Range("P35").Value = 19500
$ws.Range("S35").Value = 975

# Row 36
$ws.Range("D36").Value = 44357
$ws.Range("M36").Value = 200
$ws.Range("N36").Value = 20000
$ws.Range("O36").Value = 21000
$ws.Range("P36").Value = 20500
$ws.Range("S36").Value = 1025

# Row 37
$ws.Range("D37").Value = 44365
$ws.Range("M37").Value = 150
$ws.Range("N37").Value = 20000
$ws.Range("O37").Value = 21000
$ws.Range("P37").Value = 20500
$ws.Range("S37").Value = 1025
